# Added icons for fatal accidents and formatting changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Fatalities" column before the existing "Accident Type" column
# (old column F -> new column G).
$ws.Columns.Item(6).Insert()

# Update Weather categories for some rows: "Other" -> "Crash" / "Animal".
# New strings get interned into the shared-string table in the order they
# are first typed, so match that order here: Animal, then Crash.
$ws.Cells.Item(3, 5).Value = "Animal"
$ws.Cells.Item(6, 5).Value = "Animal"
$ws.Cells.Item(2, 5).Value = "Crash"
$ws.Cells.Item(9, 5).Value = "Crash"

# New header for column F, and set its width.
$ws.Cells.Item(1, 6).Value = "Fatalities"
$ws.Columns.Item(6).ColumnWidth = 12.33203125

# Fatalities values (0/1) for rows 2-9.
$fatalities = @(1, 0, 0, 1, 0, 0, 0, 0)
for ($i = 0; $i -lt $fatalities.Length; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $fatalities[$i]
}

# Selection moves to G5 after edits.
$ws.Range("G5").Select()
